$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7 (shifts rows 7.. down by one, copying
# formatting from the row above, matching the author's "insert row" edit).
$ws.Range("A7").EntireRow.Insert()

# Populate the new row's cells. Insert() already carried the row-6
# formatting down onto row 7 (A7/B7 get the blank "s=1" style, C7:F7 get
# the wrap-text "s=4" style), so only the text content needs to be set.
$title = "page/parcours_liste.php" + [char]10
$body = "click bouton liste parcours :" + [char]10 + "tableau de tous les parcours"
$ws.Range("C7").Value = $title + $body

$chars = $ws.Range("C7").Characters($title.Length + 1, $body.Length)
$chars.Font.Size = 8
$chars.Font.Name = "Calibri"

# Row heights: row 7 mirrors the other 3-line header rows (37.5pt); rows
# 17 and 18 (old 16/17, now shifted) lose their earlier custom heights and
# fall back to the default, while row 19 (old 18) settles on 26.25pt and
# row 20 (old 19) goes back to the default too.
$ws.Rows(7).RowHeight = 37.5
$ws.Rows("17:18").AutoFit()
$ws.Rows(19).RowHeight = 26.25
$ws.Rows(20).AutoFit()

# Match the author's final selection.
$ws.Range("B9").Select()
